# Applies the "saving residual load in excel / preparing scenario to run
# unlimited investment" edit to the AMIRIS data structure workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "times" sheet: shift the simulation window back by 28 years
#    (StartTime 2052-12-31 -> 2024-12-31, StopTime 2053-12-30 -> 2025-12-30)
# ---------------------------------------------------------------------------
$wsTimes = $wb.Worksheets.Item("times")
$wsTimes.Range("B2").Value = 45657.99861111111
$wsTimes.Range("B3").Value = 46021.99861111111

# ---------------------------------------------------------------------------
# 2. "scenario_data_emlab" sheet: target year 2053 -> 2025
# ---------------------------------------------------------------------------
$wsScenario = $wb.Worksheets.Item("scenario_data_emlab")
$wsScenario.Range("B1").Value = 2025

# ---------------------------------------------------------------------------
# 3. "conventionals" sheet: raise block size / installed power caps so the
#    model can invest without being capped by the old small blocks
# ---------------------------------------------------------------------------
$wsConv = $wb.Worksheets.Item("conventionals")
$wsConv.Range("F2").Value = 500
$wsConv.Range("G2").Value = 500
$wsConv.Range("F3").Value = 500
$wsConv.Range("G3").Value = 500

# ---------------------------------------------------------------------------
# 4. "renewables" sheet: replace the 44-row plant list with a trimmed
#    15-row list (rows 2-4 get new identifiers/capacities, rows 5-16 take
#    over the values that used to live in rows 34-45, and old rows 17-45
#    are removed entirely).
# ---------------------------------------------------------------------------
$wsRen = $wb.Worksheets.Item("renewables")

# Delete old rows 17 through 45 (29 rows) so the table ends at row 16.
$wsRen.Range("A17:I45").EntireRow.Delete()

# Row 2: identifier + capacity updated, rest stays the same (OtherPV, 0 opex)
$wsRen.Range("B2").Value = 99992100002
$wsRen.Range("C2").Value = 500

# Row 3: identifier + capacity + opex + set updated (now a WindOn plant)
$wsRen.Range("B3").Value = 99992400003
$wsRen.Range("C3").Value = 500
$wsRen.Range("D3").Value = 1.35
$wsRen.Range("E3").Value = "WindOn"

# Row 4: identifier + capacity + opex + set updated (now a WindOff plant)
$wsRen.Range("B4").Value = 99992300007
$wsRen.Range("C4").Value = 500
$wsRen.Range("D4").Value = 2.7
$wsRen.Range("E4").Value = "WindOff"

# Rows 5-16: pull up the values that used to sit in old rows 34-45
$wsRen.Range("B5").Value = 20202300034
$wsRen.Range("C5").Value = 2591.333333333333

$wsRen.Range("B6").Value = 20152100032
$wsRen.Range("C6").Value = 18148.27119466832

$wsRen.Range("B7").Value = 20152300053
$wsRen.Range("C7").Value = 2591.333333333333
$wsRen.Range("D7").Value = 2.7
$wsRen.Range("E7").Value = "WindOff"

$wsRen.Range("B8").Value = 20112400035
$wsRen.Range("C8").Value = 17185.46324999998
$wsRen.Range("D8").Value = 1.35
$wsRen.Range("E8").Value = "WindOn"

$wsRen.Range("B9").Value = 20102300046
$wsRen.Range("C9").Value = 2591.333333333333
$wsRen.Range("D9").Value = 2.7
$wsRen.Range("E9").Value = "WindOff"

$wsRen.Range("B10").Value = 20102100052
$wsRen.Range("C10").Value = 18148.27119466832

$wsRen.Range("B11").Value = 20062400023
$wsRen.Range("C11").Value = 17185.46324999998
$wsRen.Range("D11").Value = 1.35
$wsRen.Range("E11").Value = "WindOn"

$wsRen.Range("B12").Value = 20052100045
$wsRen.Range("C12").Value = 18148.27119466832

$wsRen.Range("B13").Value = 20012400047
$wsRen.Range("C13").Value = 17185.46324999998
$wsRen.Range("D13").Value = 1.35
$wsRen.Range("E13").Value = "WindOn"

$wsRen.Range("B14").Value = 19691200028
$wsRen.Range("C14").Value = 2940.25
$wsRen.Range("E14").Value = "RunOfRiver"

$wsRen.Range("B15").Value = 19641200048
$wsRen.Range("C15").Value = 2940.25
$wsRen.Range("E15").Value = "RunOfRiver"

$wsRen.Range("B16").Value = 19591200040
$wsRen.Range("C16").Value = 2940.25
$wsRen.Range("E16").Value = "RunOfRiver"

# ---------------------------------------------------------------------------
# 5. "storages" sheet: raise installed power cap 100 -> 500
# ---------------------------------------------------------------------------
$wsStor = $wb.Worksheets.Item("storages")
$wsStor.Range("H2").Value = 500

# ---------------------------------------------------------------------------
# 6. "biogas" sheet: raise installed power cap 100 -> 500
# ---------------------------------------------------------------------------
$wsBio = $wb.Worksheets.Item("biogas")
$wsBio.Range("C2").Value = 500
